# Apply cryptos.xlsx data refresh (GitHub Actions update)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.058.48"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.53%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.675.42"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.40%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.18"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.20%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.517"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.01%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.97%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "21.42"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +5.62%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0622"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0883"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.60%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.913.26"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.49%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.677.77"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.26%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.12"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.13%  "
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.97%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.24"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.06%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.043.88"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.43%  "
$ws.Range("B18").Value = "BitcoinCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "237.98"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.88%  "
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.15"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +2.17%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0737"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.56%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.47"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.41%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.92%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.13"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.37%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.73"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.86%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +2.76%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.81%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.19%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.47%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.24%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.35"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.62%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.532.84"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +5.55%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.18"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.65%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +2.88%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.38"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.11%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.592"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.13%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.918"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +2.51%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +2.13%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.07"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +2.82%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.02%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "67.78"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +2.23%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.26"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.39%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -3.92%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.820.57"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.45%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.782"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.32%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "90.72"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.24%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.98%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +2.25%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.05"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +5.80%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.41%  "
